# 11 May - Noche
# Reassign "Materia" (column E) and "Docente" (column F) pairs among the
# rows that belong to the same student block on the "Blancos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blancos")

# Map: row number -> @(new Materia text, new Docente text)
$changes = @{
    4  = @("CÁLCULO DIFERENCIAL", "Rodríguez Román Leticia")
    5  = @("ECOLOGÍA", "Rivera Cruz Ezequiel")
    6  = @("FÍSICA I", "González Sánchez Rene Aurelio")
    8  = @("MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "Velasco Sanchez David")
    9  = @("CÁLCULO DIFERENCIAL", "Rodríguez Román Leticia")
    13 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    14 = @("MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "Velasco Sanchez David")
    20 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    21 = @("CÁLCULO DIFERENCIAL", "Rodríguez Román Leticia")
    22 = @("CÁLCULO DIFERENCIAL", "Rodríguez Román Leticia")
    23 = @("MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "Velasco Sanchez David")
    24 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    26 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    27 = @("CÁLCULO DIFERENCIAL", "Rodríguez Román Leticia")
    28 = @("ECOLOGÍA", "Rivera Cruz Ezequiel")
    29 = @("MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "Velasco Sanchez David")
    30 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    31 = @("CÁLCULO DIFERENCIAL", "Rodríguez Román Leticia")
    40 = @("MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "Velasco Sanchez David")
    41 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    44 = @("ECOLOGÍA", "Rivera Cruz Ezequiel")
    45 = @("CÁLCULO DIFERENCIAL", "Rodríguez Román Leticia")
    46 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    48 = @("CÁLCULO DIFERENCIAL", "Rodríguez Román Leticia")
    49 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    50 = @("CÁLCULO DIFERENCIAL", "Rodríguez Román Leticia")
    55 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    56 = @("MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "Velasco Sanchez David")
    59 = @("MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "Velasco Sanchez David")
    60 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    63 = @("CÁLCULO DIFERENCIAL", "Rodríguez Román Leticia")
    64 = @("MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTROMAGNÉTICO", "Velasco Sanchez David")
    65 = @("FÍSICA I", "González Sánchez Rene Aurelio")
}

foreach ($row in $changes.Keys) {
    $pair = $changes[$row]
    $ws.Cells.Item($row, 5).Value = $pair[0]
    $ws.Cells.Item($row, 6).Value = $pair[1]
}
